$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New summary rows (14-17) -------------------------------------------
# Row 14: Average of SW(S*)/SW(OPT)
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

# Row 15: Average of SC(S*)/SC(OPT)
$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

# Row 16: Worst of SW(S*)/SW(OPT)
$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

# Row 17: Worst of SC(S*)/SC(OPT)
$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Bold, size 12, vertically centered font for the new summary values, with a
# taller row height to match the bigger font.
$summary = $ws.Range("B14:B17")
$summary.VerticalAlignment = -4108
$summary.Font.Bold = $true
$summary.Font.Size = 12
$summary.EntireRow.RowHeight = 15.6

# --- Row 12: average of column J (|S*|/n) across the data rows ---------
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"
$ws.Range("J12").Font.Bold = $true

# Selection matches the author's last interaction
$ws.Range("A14:B17").Select() | Out-Null

# Page setup: portrait orientation, paper size 9 (A4)
$ws.PageSetup.Orientation = 1
$ws.PageSetup.PaperSize = 9
